# Append new Lancers listing captured at 2025-09-21 12:32:58 JST.
# A new row is inserted right after the header/first listing (becomes row 3),
# every row below shifts down by one, and all "fetched at" timestamps in
# column A move to the new run time. A new hyperlink is appended for the
# newly last row (old last row, now row 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-09-21 12:32:58"

# Clear existing hyperlinks up front -- they are keyed off cell position and
# would otherwise end up pointing at the wrong row once the data shifts.
$ws.Range("A1").Hyperlinks.Delete()

# Widen the price column slightly to fit the new (longer) price range text.
$ws.Columns.Item(4).ColumnWidth = 29.14

# Final data for rows 2..7 (row 1 is the header and is untouched).
$rows = @(
    @($timestamp, "GPTsを使ってAi個人会話webアプリの開発", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5397594", 478, "🔥AI,GPT ◆開発 ◇アプリ"),
    @($timestamp, "【時給1万円超】AI活用を伝える講師パートナー募集", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5397680", 310, "🔥AI,Ai"),
    @($timestamp, "【急募】検査報告書自動発行ツール開発のエキスパート募集", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5397543", 123, "◆ツール,開発"),
    @($timestamp, "【急募】WordPress開発者を探しています!魅力的な案件です", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5397452", 88, "◆開発 ○WordPress"),
    @($timestamp, "初回 電子秤からのデータ抽出のシステム化", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5397615", 33, ""),
    @($timestamp, "【急募】JotformとGoogleスプレッドシート連携のエキスパート募集!", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5395809", 10, "")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value2 = $data[0]
    $ws.Cells.Item($r, 2).Value2 = $data[1]
    $ws.Cells.Item($r, 3).Value2 = $data[2]
    $ws.Cells.Item($r, 4).Value2 = $data[3]
    $ws.Cells.Item($r, 5).Value2 = $data[4]
    $ws.Cells.Item($r, 6).Value2 = $data[5]
    $ws.Cells.Item($r, 7).Value2 = $data[6]
    # Two of the six rows have no "skill summary" text at all (no H cell in
    # the source data) -- writing "" leaves the cell empty/absent, matching.
    $ws.Cells.Item($r, 8).Value2 = $data[7]

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $data[5])
    # `Hyperlinks.Add` mints a fresh (duplicate) cell-style entry; pin the
    # cell back to the shared built-in "Hyperlink" style so the whole
    # column keeps using the same style index as before the edit.
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
